$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.062.91'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.313.43'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.73%  '
$ws.Range('E7').Value = '  +2.82%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.523'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '2.673.60'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '2.294.40'
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('E17').Value = '  -3.03%  '
$ws.Range('D18').Value = '42.978.82'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.34%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '168.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('E30').Value = '  -12.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.45%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.36'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.65%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0694'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('E40').Value = '  +0.97%  '
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('D42').Value = '2.001.73'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.54'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('D49').Value = '2.539.55'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.30%  '
$ws.Range('E51').Value = '  +0.34%  '
